# Refined metadata to be additional tab
# 1) Update the time_taken timestamps on the existing "data" sheet.
# 2) Add a new "metadata" worksheet (after "data") with panel query metadata.

$wb = $excel.ActiveWorkbook
$dataSheet = $wb.Worksheets.Item("data")

# --- Update time_taken column (F) on the "data" sheet ---
$dataSheet.Range("F2").Value = "2021-10-05 14:33:42.845621"
$dataSheet.Range("F3").Value = "2021-10-05 14:33:42.845629"
$dataSheet.Range("F4").Value = "2021-10-05 14:33:42.845632"

# --- Add the new "metadata" worksheet right after "data" ---
$ws = $wb.Worksheets.Add($null, $dataSheet)
$ws.Name = "metadata"

# Match the page margins used on the "data" sheet (values are in points).
$ws.PageSetup.LeftMargin = 0.75 * 72
$ws.PageSetup.RightMargin = 0.75 * 72
$ws.PageSetup.TopMargin = 1 * 72
$ws.PageSetup.BottomMargin = 1 * 72
$ws.PageSetup.HeaderMargin = 0.5 * 72
$ws.PageSetup.FooterMargin = 0.5 * 72

# Header row (row 1), bold + bordered + centered like the header row in "data"
$headerCells = @("B1","C1","D1","E1","F1","G1")
$headerValues = @("data_name","data_id","data_version","data_version_created","panel_query_time","panel_get_request")
for ($i = 0; $i -lt $headerCells.Length; $i++) {
    $cell = $ws.Range($headerCells[$i])
    $cell.Value = $headerValues[$i]
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
    $cell.Borders.LineStyle = 1
    $cell.Font.Bold = $true
}

# Data row (row 2)
$a2 = $ws.Range("A2")
$a2.Value = 0
$a2.HorizontalAlignment = -4108
$a2.VerticalAlignment = -4160
$a2.Borders.LineStyle = 1
$a2.Font.Bold = $true

$ws.Range("B2").Value = "Diabetes Insipidus"
$ws.Range("C2").Value = 3445
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "1.1"
$ws.Range("E2").Value = "2021-03-18T05:04:20.288275Z"
$ws.Range("F2").Value = "2021-10-05 14:33:42.841753"
$ws.Range("G2").Value = "https://panelapp.agha.umccr.org/api/v1/panels/3445/?format=json"

# Keep "data" as the active/selected sheet, same as before the edit.
$dataSheet.Activate()
